$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge the two existing column groups (A and B:D) into a single
# A:D column band with a new width, same as the source edit.
$ws.Columns("A:D").ColumnWidth = 27.59

# Add the new row 8 with the translated label repeated across A:D,
# wrapped and taller to fit the two-line string.
$ws.Range("A8:D8").Value = "View/Edit Region Business Rules`n"
$ws.Range("A8:D8").WrapText = $true
$ws.Range("A8:D8").RowHeight = 24
